$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "A Lag" row (row 2): C coefficient lag value
$ws.Range("B2").Value = "-0.37***"

# Update the "C Lag" row (row 3): A coefficient lag value and C coefficient lag value
$ws.Range("B3").Value = "-3.46***"
$ws.Range("C3").Value = "-0.81***"
